$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$win = $excel.ActiveWindow

# Reworked style-setting for pupil-course tables:
# 1) Freeze panes at column D / row 7 boundary (so that column E and row 8
#    become the first scrollable column/row), matching the header block.
$ws.Range("A1").Select()
$ws.Range("E8").Select()
$win.FreezePanes = $true

# Put the final active selection on M8 in the bottom-right (scrollable) pane,
# as in the target view state.
$ws.Range("M8").Select()

# 2) Clear the redundant "MFN" subject-style values from M12:M15 while
#    keeping their existing cell style (s="8") intact.
$ws.Range("M12:M15").ClearContents()
